# Update the "想去人数" (attendance) numbers across the workbook's sheets
# to reflect the latest output snapshot (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 200
$ws1.Range("F3").Value = 534
$ws1.Range("F4").Value = 43
$ws1.Range("F9").Value = 343
$ws1.Range("F10").Value = 3361
$ws1.Range("F11").Value = 37

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 93

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 200
$ws4.Range("F3").Value = 93
$ws4.Range("F4").Value = 534
$ws4.Range("F5").Value = 43
$ws4.Range("F10").Value = 343
$ws4.Range("F11").Value = 3361
$ws4.Range("F12").Value = 37
